$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 650
$ws.Range("I88").Value = 650
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 650
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -244
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 650
$ws.Range("I91").Value = 650
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 650
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = 754
$ws.Range("N91").ClearContents()
$ws.Range("H103").Value = 564.2174
$ws.Range("I103").Value = 352.5625
$ws.Range("J103").Value = 1048
$ws.Range("K103").Value = 1057.6875
$ws.Range("L103").Value = 3144
$ws.Range("M103").Value = -471.6875
$ws.Range("N103").Value = -4316
$ws.Range("H132").Value = 3098.6
$ws.Range("I132").Value = 2824.375
$ws.Range("J132").Value = 4195.5
$ws.Range("K132").Value = 8473.125
$ws.Range("L132").Value = 12586.5
$ws.Range("M132").Value = -5943.125
$ws.Range("N132").Value = -17646.5
$ws.Range("H137").Value = 37871.566
$ws.Range("I137").Value = 45169.95
$ws.Range("J137").Value = 3204.25
$ws.Range("K137").Value = 135509.85
$ws.Range("L137").Value = 9612.75
$ws.Range("M137").Value = -132959.85
$ws.Range("N137").Value = -14712.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 1199.6666
$ws.Range("J13").Value = 1199.6666
$ws.Range("L13").Value = 1199.6666
$ws.Range("N13").Value = -1487.6666
$ws.Range("H32").Value = 17097650
$ws.Range("I32").Value = 17634582
$ws.Range("K32").Value = 17634582
$ws.Range("M32").Value = -17634295
$ws.Range("H61").Value = 3282.5386
$ws.Range("I61").Value = 2669.5833
$ws.Range("K61").Value = 2669.5833
$ws.Range("M61").Value = -2457.5833
$ws.Range("H123").Value = 74995
$ws.Range("J123").Value = 74995
$ws.Range("L123").Value = 74995
$ws.Range("N123").Value = -84795
$ws.Range("H136").Value = 3282.5386
$ws.Range("I136").Value = 2669.5833
$ws.Range("K136").Value = 8008.749899999999
$ws.Range("M136").Value = -5458.749899999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2660.9333
$ws.Range("I20").Value = 2376.375
$ws.Range("K20").Value = 2376.375
$ws.Range("M20").Value = -2129.375
$ws.Range("H122").Value = 115500
$ws.Range("J122").Value = 115500
$ws.Range("L122").Value = 115500
$ws.Range("N122").Value = -125300

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 2000
$ws.Range("J38").Value = 2000
$ws.Range("L38").Value = 2000
$ws.Range("N38").Value = -2754
$ws.Range("H46").Value = 2000
$ws.Range("J46").Value = 2000
$ws.Range("L46").Value = 2000
$ws.Range("N46").Value = -2422
$ws.Range("H62").Value = 4737.2
$ws.Range("I62").Value = 4499.5
$ws.Range("J62").Value = 4895.6665
$ws.Range("K62").Value = 4499.5
$ws.Range("L62").Value = 4895.6665
$ws.Range("M62").Value = -3875.5
$ws.Range("N62").Value = -6143.6665
$ws.Range("H65").Value = 4737.2
$ws.Range("I65").Value = 4499.5
$ws.Range("J65").Value = 4895.6665
$ws.Range("K65").Value = 22497.5
$ws.Range("L65").Value = 24478.3325
$ws.Range("M65").Value = -19377.5
$ws.Range("N65").Value = -30718.3325
$ws.Range("H68").Value = 70000
$ws.Range("J68").Value = 70000
$ws.Range("L68").Value = 70000
$ws.Range("N68").Value = -71498
$ws.Range("H71").Value = 70000
$ws.Range("J71").Value = 70000
$ws.Range("L71").Value = 210000
$ws.Range("N71").Value = -217488
$ws.Range("H116").Value = 167999
$ws.Range("J116").Value = 167999
$ws.Range("L116").Value = 167999
$ws.Range("N116").Value = -177177
$ws.Range("H134").Value = 2978
$ws.Range("I134").Value = 2000
$ws.Range("K134").Value = 6000
$ws.Range("M134").Value = -3465

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 2000
$ws.Range("I82").Value = 2000
$ws.Range("K82").Value = 6000
$ws.Range("M82").Value = -5594
$ws.Range("H85").Value = 2000
$ws.Range("I85").Value = 2000
$ws.Range("K85").Value = 6000
$ws.Range("M85").Value = -4596
$ws.Range("H93").Value = 1337.5
$ws.Range("H131").Value = 1829.1428
$ws.Range("J131").Value = 1867.625
$ws.Range("L131").Value = 5602.875
$ws.Range("N131").Value = -15682.875

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 6002249
$ws.Range("I11").Value = 8574213
$ws.Range("K11").Value = 8574213
$ws.Range("M11").Value = -8574074
$ws.Range("H135").Value = 215000.5
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 215000.5
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 215000.5
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -225140.5
$ws.Range("H139").Value = 217499.83
$ws.Range("J139").Value = 217499.83
$ws.Range("L139").Value = 217499.83
$ws.Range("N139").Value = -227779.83

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7317.4
$ws.Range("I7").Value = 7510.4546
$ws.Range("J7").Value = 7081.4443
$ws.Range("K7").Value = 7510.4546
$ws.Range("L7").Value = 7081.4443
$ws.Range("M7").Value = -7398.4546
$ws.Range("N7").Value = -7305.4443
$ws.Range("H22").Value = 3011.889
$ws.Range("J22").Value = 4300.4
$ws.Range("L22").Value = 4300.4
$ws.Range("N22").Value = -4890.4
$ws.Range("H27").Value = 3011.889
$ws.Range("J27").Value = 4300.4
$ws.Range("L27").Value = 4300.4
$ws.Range("N27").Value = -4514.4
$ws.Range("H32").Value = 27666.334
$ws.Range("I32").Value = 27666.334
$ws.Range("K32").Value = 27666.334
$ws.Range("M32").Value = -27349.334
$ws.Range("H40").Value = 4665.6665
$ws.Range("I40").Value = 4098.8
$ws.Range("J40").Value = 7500
$ws.Range("K40").Value = 4098.8
$ws.Range("L40").Value = 7500
$ws.Range("M40").Value = -3962.8
$ws.Range("N40").Value = -7772
$ws.Range("H46").Value = 8568.5
$ws.Range("I46").Value = 3957
$ws.Range("J46").Value = 8987.727999999999
$ws.Range("K46").Value = 3957
$ws.Range("L46").Value = 8987.727999999999
$ws.Range("M46").Value = -3769
$ws.Range("N46").Value = -9363.727999999999
$ws.Range("H68").Value = 2498.5
$ws.Range("J68").Value = 2668.3333
$ws.Range("L68").Value = 2668.3333
$ws.Range("N68").Value = -4166.3333
$ws.Range("H71").Value = 2498.5
$ws.Range("J71").Value = 2668.3333
$ws.Range("L71").Value = 13341.6665
$ws.Range("N71").Value = -20829.6665
$ws.Range("H93").Value = 1508.6818
$ws.Range("I93").Value = 1418.75
$ws.Range("J93").Value = 1616.6
$ws.Range("K93").Value = 1418.75
$ws.Range("L93").Value = 1616.6
$ws.Range("M93").Value = -170.75
$ws.Range("N93").Value = -4112.6
$ws.Range("H126").Value = 7317.4
$ws.Range("I126").Value = 7510.4546
$ws.Range("J126").Value = 7081.4443
$ws.Range("K126").Value = 22531.3638
$ws.Range("L126").Value = 21244.3329
$ws.Range("M126").Value = -20061.3638
$ws.Range("N126").Value = -26184.3329

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 250010500
$ws.Range("J14").Value = 10999
$ws.Range("L14").Value = 10999
$ws.Range("N14").Value = -11335
$ws.Range("H43").Value = 105490
$ws.Range("I43").Value = 105490
$ws.Range("K43").Value = 105490
$ws.Range("M43").Value = -105341
$ws.Range("H49").Value = 8056
$ws.Range("I49").Value = 8056
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 8056
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -7826
$ws.Range("N49").ClearContents()
$ws.Range("H136").Value = 930.6
$ws.Range("I136").Value = 930.6
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2791.8
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -241.8000000000002
$ws.Range("N136").ClearContents()

